$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value
$ws.Range("B2").Value = 120339

# Add new row data: A3 = 2, B3 = 120340
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 120340

# Update the selection to B2:B3 with active cell B2
$ws.Range("B2:B3").Select()
